# Update provincias_spain data: refreshed case counts causing re-sort of several
# rows (by "Casos totales" descending) plus a refreshed "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / "last updated" banner (A1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 6 de Abril de 2020 a las 12:52"

# Row 8 - Ciudad Real (stays in place, values refreshed)
$ws.Cells.Item(8,2).Value = 4125
$ws.Cells.Item(8,3).Value = 1259
$ws.Cells.Item(8,4).Value = 8211
$ws.Cells.Item(8,5).Value = 364

# Rows 15-17 - Albacete/La Rioja/Alacant swap order due to re-sort
$ws.Cells.Item(15,1).Value = "Albacete"
$ws.Cells.Item(15,2).Value = 2751
$ws.Cells.Item(15,3).Value = 1259
$ws.Cells.Item(15,4).Value = 8211
$ws.Cells.Item(15,5).Value = 241

$ws.Cells.Item(16,1).Value = "La Rioja"
$ws.Cells.Item(16,2).Value = 2719
$ws.Cells.Item(16,3).Value = 964
$ws.Cells.Item(16,4).Value = 1614
$ws.Cells.Item(16,5).Value = 141

$ws.Cells.Item(17,1).Value = "Alacant/Alicante"
$ws.Cells.Item(17,2).Value = 2673
$ws.Cells.Item(17,3).Value = 284
$ws.Cells.Item(17,4).Value = 2120
$ws.Cells.Item(17,5).Value = 269

# Row 19 - Zaragoza (stays in place, values refreshed)
$ws.Cells.Item(19,2).Value = 2520
$ws.Cells.Item(19,3).Value = 408
$ws.Cells.Item(19,4).Value = 1895
$ws.Cells.Item(19,5).Value = 217

# Row 20 - Toledo (stays in place, values refreshed)
$ws.Cells.Item(20,2).Value = 2283
$ws.Cells.Item(20,3).Value = 1259
$ws.Cells.Item(20,4).Value = 8211
$ws.Cells.Item(20,5).Value = 307

# Row 31 - Caceres (stays in place, values refreshed)
$ws.Cells.Item(31,2).Value = 1379
$ws.Cells.Item(31,3).Value = 87
$ws.Cells.Item(31,4).Value = 1097
$ws.Cells.Item(31,5).Value = 195

# Rows 38-39 - Guadalajara/Castello swap order due to re-sort
$ws.Cells.Item(38,1).Value = "Guadalajara"
$ws.Cells.Item(38,2).Value = 873
$ws.Cells.Item(38,3).Value = 1259
$ws.Cells.Item(38,4).Value = 8211
$ws.Cells.Item(38,5).Value = 124

$ws.Cells.Item(39,1).Value = "Castello/Castellon"
$ws.Cells.Item(39,2).Value = 869
$ws.Cells.Item(39,3).Value = 85
$ws.Cells.Item(39,4).Value = 714
$ws.Cells.Item(39,5).Value = 70

# Rows 42-43 - Badajoz/Avila swap order due to re-sort
$ws.Cells.Item(42,1).Value = "Badajoz"
$ws.Cells.Item(42,2).Value = 689
$ws.Cells.Item(42,3).Value = 136
$ws.Cells.Item(42,4).Value = 520
$ws.Cells.Item(42,5).Value = 33

$ws.Cells.Item(43,1).Value = "Avila"
$ws.Cells.Item(43,2).Value = 679
$ws.Cells.Item(43,3).Value = 214
$ws.Cells.Item(43,4).Value = 392
$ws.Cells.Item(43,5).Value = 73

# Rows 46-47 - Cuenca/Huelva swap order due to re-sort
$ws.Cells.Item(46,1).Value = "Cuenca"
$ws.Cells.Item(46,2).Value = 570
$ws.Cells.Item(46,3).Value = 1259
$ws.Cells.Item(46,4).Value = 8211
$ws.Cells.Item(46,5).Value = 96

$ws.Cells.Item(47,1).Value = "Huelva"
$ws.Cells.Item(47,2).Value = 553
$ws.Cells.Item(47,3).Value = 18
$ws.Cells.Item(47,4).Value = 524
$ws.Cells.Item(47,5).Value = 11

# Row 50 - Huesca (stays in place, values refreshed)
$ws.Cells.Item(50,2).Value = 417
$ws.Cells.Item(50,3).Value = 57
$ws.Cells.Item(50,4).Value = 323
$ws.Cells.Item(50,5).Value = 37

# Row 51 - Teruel (stays in place, values refreshed)
$ws.Cells.Item(51,2).Value = 380
$ws.Cells.Item(51,3).Value = 93
$ws.Cells.Item(51,4).Value = 258
